$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E8").Value = 10
$ws.Range("E9").Value = 10
$ws.Range("E10").Value = 10
$ws.Range("E11").Value = 40

$ws.Range("E12").Select()
